$wb = $excel.ActiveWorkbook

# "babies" sheet: add a "treatment" header column in D1
$wsBabies = $wb.Worksheets.Item("babies")
$wsBabies.Range("D1").Value = "treatment"

# "rooms" sheet: add a "treatment" header column in I1
$wsRooms = $wb.Worksheets.Item("rooms")
$wsRooms.Range("I1").Value = "treatment"

# Update selections to reflect new active cells
$wsBabies.Range("D1").Select()
$wsRooms.Range("I1").Select()

# Make "babies" the active sheet (first sheet)
$wsBabies.Activate()
